$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (set in G, F, E order so the shared-string table is
# built in the same order as the target workbook: "l27", "l25.0", "l23.0")
$ws.Cells.Item(1, 7).Value = "l27"
$ws.Cells.Item(1, 6).Value = "l25.0"
$ws.Cells.Item(1, 5).Value = "l23.0"

# New data columns E, F, G mirror the existing B, C, D values
$ws.Cells.Item(2, 5).Value = 0.046
$ws.Cells.Item(2, 6).Value = 0.051
$ws.Cells.Item(2, 7).Value = 0.057

$ws.Cells.Item(3, 5).Value = 0.077
$ws.Cells.Item(3, 6).Value = 0.085
$ws.Cells.Item(3, 7).Value = 0.094

$ws.Cells.Item(4, 5).Value = 0.11
$ws.Cells.Item(4, 6).Value = 0.122
$ws.Cells.Item(4, 7).Value = 0.134

$ws.Cells.Item(5, 5).Value = 0.145
$ws.Cells.Item(5, 6).Value = 0.159
$ws.Cells.Item(5, 7).Value = 0.174

$ws.Cells.Item(6, 5).Value = 0.176
$ws.Cells.Item(6, 6).Value = 0.194
$ws.Cells.Item(6, 7).Value = 0.211

$ws.Cells.Item(7, 5).Value = 0.205
$ws.Cells.Item(7, 6).Value = 0.224
$ws.Cells.Item(7, 7).Value = 0.244

$ws.Cells.Item(8, 5).Value = 0.229
$ws.Cells.Item(8, 6).Value = 0.25
$ws.Cells.Item(8, 7).Value = 0.272

$ws.Cells.Item(9, 5).Value = 0.248
$ws.Cells.Item(9, 6).Value = 0.271
$ws.Cells.Item(9, 7).Value = 0.295

$ws.Cells.Item(10, 5).Value = 0.264
$ws.Cells.Item(10, 6).Value = 0.288
$ws.Cells.Item(10, 7).Value = 0.312

$ws.Cells.Item(11, 5).Value = 0.276
$ws.Cells.Item(11, 6).Value = 0.3
$ws.Cells.Item(11, 7).Value = 0.325

$ws.Cells.Item(12, 5).Value = 0.284
$ws.Cells.Item(12, 6).Value = 0.309
$ws.Cells.Item(12, 7).Value = 0.335

# Move the active selection to E2, matching the saved view state
$ws.Range("E2").Select() | Out-Null
